$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9678354263305664
$ws.Range("B1").Value = 3.046417713165283
$ws.Range("C1").Value = 4.128415584564209
$ws.Range("D1").Value = 2.059024810791016
$ws.Range("E1").Value = 1.225955247879028
